# "Correction optimizer 1000 runs data"
#
# Both "Optimizer | 1000 run" sheets (Ascended and Descend) had the H and I
# columns (the cached gas-consumption simulation figures) corrected. Column
# H (a constant baseline figure) drops from 2599085 to 2590619 on every
# data row, and column I (the per-row variable figure) is updated to the
# newly-corrected value. Columns J/L/M/R recompute automatically from these
# via formulas, and the chart caches that plot R (and mirror it on the
# "Graphs" sheet) refresh automatically on recalculation.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Optimizer | 1000 run - Ascended" (rows 3-13): new H / I values
# ---------------------------------------------------------------------
$wsAsc = $wb.Worksheets.Item("Optimizer | 1000 run - Ascended")

$data1000Asc = @(
    @(3, 62268),
    @(4, 180609),
    @(5, 342594),
    @(6, 536679),
    @(7, 762864),
    @(8, 1021149),
    @(9, 1311534),
    @(10, 1634019),
    @(11, 1988604),
    @(12, 2375289),
    @(13, 2789274)
)

# H3:H13 is a single constant value across every row of this table
$wsAsc.Range("H3:H13").Value = 2590619

foreach ($pair in $data1000Asc) {
    $r = $pair[0]
    $i = $pair[1]
    $wsAsc.Range("I$r").Value = $i
}

# ---------------------------------------------------------------------
# 2. "Optimizer | 1000 run - Descend" (rows 3-23): new H / I values
# ---------------------------------------------------------------------
$wsDesc = $wb.Worksheets.Item("Optimizer | 1000 run - Descend")

$data1000Desc = @(
    @(3, 62317),
    @(4, 180091),
    @(5, 339356),
    @(6, 528521),
    @(7, 724334),
    @(8, 747586),
    @(9, 771137),
    @(10, 794987),
    @(11, 819136),
    @(12, 843584),
    @(13, 868331),
    @(14, 893377),
    @(15, 918722),
    @(16, 944366),
    @(17, 970309),
    @(18, 996551),
    @(19, 1275416),
    @(20, 1584181),
    @(21, 1922846),
    @(22, 2291411),
    @(23, 2685076)
)

$wsDesc.Range("H3:H23").Value = 2590619

foreach ($pair in $data1000Desc) {
    $r = $pair[0]
    $i = $pair[1]
    $wsDesc.Range("I$r").Value = $i
}

# ---------------------------------------------------------------------
# 3. Update the selections left behind on each of the touched sheets,
#    and make "Optimizer | 1000 run - Descend" the active tab (it was
#    "Graphs" before).
# ---------------------------------------------------------------------
$wsAsc.Activate()
$wsAsc.Range("I16").Select()

$wsDesc.Activate()
$wsDesc.Range("K19").Select()
